$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Epoch2")
$ws.Rows("2:16").Insert()
$ws.Range("A2:A16").ClearFormats()
$ws.Range("A2").Value = 'Epoch: 77, Val_loss:  2.33, time:104.85, test_Acc:  26.99, test_bacc:  27.81, test_f1:  18.68'
$ws.Range("A3").Value = 'Epoch: 76, Val_loss:  2.33, time:105.00, test_Acc:  26.99, test_bacc:  27.81, test_f1:  18.68'
$ws.Range("A4").Value = 'Epoch: 75, Val_loss:  2.31, time:105.35, test_Acc:  26.99, test_bacc:  27.81, test_f1:  18.68'
$ws.Range("A5").Value = 'Epoch: 74, Val_loss:  2.31, time:104.89, test_Acc:  26.99, test_bacc:  27.81, test_f1:  18.68'
$ws.Range("A6").Value = 'Epoch: 73, Val_loss:  2.29, time:105.24, test_Acc:  26.99, test_bacc:  27.81, test_f1:  18.68'
$ws.Range("A7").Value = 'Epoch: 72, Val_loss:  2.30, time:104.71, test_Acc:  26.99, test_bacc:  27.81, test_f1:  18.68'
$ws.Range("A8").Value = 'Epoch: 71, Val_loss:  2.28, time:104.76, test_Acc:  26.99, test_bacc:  27.81, test_f1:  18.68'
$ws.Range("A9").Value = 'Epoch: 70, Val_loss:  2.28, time:104.94, test_Acc:  26.99, test_bacc:  27.81, test_f1:  18.68'
$ws.Range("A10").Value = 'Epoch: 69, Val_loss:  2.26, time:104.56, test_Acc:  26.99, test_bacc:  27.81, test_f1:  18.68'
$ws.Range("A11").Value = 'Epoch: 68, Val_loss:  2.26, time:104.56, test_Acc:  26.99, test_bacc:  27.81, test_f1:  18.68'
$ws.Range("A12").Value = 'Epoch: 67, Val_loss:  2.24, time:105.64, test_Acc:  26.99, test_bacc:  27.81, test_f1:  18.68'
$ws.Range("A13").Value = 'Epoch: 66, Val_loss:  2.24, time:104.68, test_Acc:  26.99, test_bacc:  27.81, test_f1:  18.68'
$ws.Range("A14").Value = 'Epoch: 65, Val_loss:  2.22, time:104.34, test_Acc:  26.99, test_bacc:  27.81, test_f1:  18.68'
$ws.Range("A15").Value = 'Epoch: 64, Val_loss:  2.22, time:106.29, test_Acc:  26.99, test_bacc:  27.81, test_f1:  18.68'
$ws.Range("A16").Value = 'Epoch: 63, Val_loss:  2.20, time:105.29, test_Acc:  26.99, test_bacc:  27.81, test_f1:  18.68'
